$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build a 2D array holding the new rows of stock/NAV data to append
$arr = New-Object 'object[,]' 20,10
$arr[0,0] = "2024-09-02"
$arr[0,2] = [double]"1092.650024414062"
$arr[0,3] = [double]"715.0499877929688"
$arr[0,4] = [double]"251.3500061035156"
$arr[0,5] = [double]"490.5"
$arr[0,6] = [double]"1505.25"
$arr[0,7] = [double]"29322.30020141602"
$arr[0,8] = [double]"0"
$arr[0,9] = [double]"240.2363002559728"

$arr[1,0] = "2024-09-03"
$arr[1,2] = [double]"1085.099975585938"
$arr[1,3] = [double]"710.7999877929688"
$arr[1,4] = [double]"251"
$arr[1,5] = [double]"488.8500061035156"
$arr[1,6] = [double]"1509"
$arr[1,7] = [double]"29201.39971923828"
$arr[1,8] = [double]"-0.004123158188384413"
$arr[1,9] = [double]"239.2457679874252"

$arr[2,0] = "2024-09-04"
$arr[2,2] = [double]"1080.449951171875"
$arr[2,3] = [double]"722.4000244140625"
$arr[2,4] = [double]"250.5"
$arr[2,5] = [double]"484.1499938964844"
$arr[2,6] = [double]"1488.099975585938"
$arr[2,7] = [double]"29200.89978027344"
$arr[2,8] = [double]"-1.712037675078922e-05"
$arr[2,9] = [double]"239.2416720097412"

$arr[3,0] = "2024-09-05"
$arr[3,2] = [double]"1069.150024414062"
$arr[3,3] = [double]"733.8499755859375"
$arr[3,4] = [double]"251.1499938964844"
$arr[3,5] = [double]"495.6499938964844"
$arr[3,6] = [double]"1447.599975585938"
$arr[3,7] = [double]"29262.3996887207"
$arr[3,8] = [double]"0.002106096350113556"
$arr[3,9] = [double]"239.745538021956"

$arr[4,0] = "2024-09-06"
$arr[4,2] = [double]"1049.349975585938"
$arr[4,3] = [double]"718.9000244140625"
$arr[4,4] = [double]"247.8000030517578"
$arr[4,5] = [double]"483"
$arr[4,6] = [double]"1418.050048828125"
$arr[4,7] = [double]"28702.20024108887"
$arr[4,8] = [double]"-0.01914400232349252"
$arr[4,9] = [double]"235.1558488850167"

$arr[5,0] = "2024-09-09"
$arr[5,2] = [double]"1038.699951171875"
$arr[5,3] = [double]"700.1500244140625"
$arr[5,4] = [double]"243.8999938964844"
$arr[5,5] = [double]"474.75"
$arr[5,6] = [double]"1411.849975585938"
$arr[5,7] = [double]"28242.6496887207"
$arr[5,8] = [double]"-0.01601098691069303"
$arr[5,9] = [double]"231.3907716665458"

$arr[6,0] = "2024-09-10"
$arr[6,2] = [double]"1035.800048828125"
$arr[6,3] = [double]"713.4000244140625"
$arr[6,4] = [double]"248.25"
$arr[6,5] = [double]"478.7999877929688"
$arr[6,6] = [double]"1424.449951171875"
$arr[6,7] = [double]"28522.85046386719"
$arr[6,8] = [double]"0.009921192885042528"
$arr[6,9] = [double]"233.6864441440684"

$arr[7,0] = "2024-09-11"
$arr[7,2] = [double]"976.2999877929688"
$arr[7,3] = [double]"725.4000244140625"
$arr[7,4] = [double]"241.5500030517578"
$arr[7,5] = [double]"472.2000122070312"
$arr[7,6] = [double]"1399.599975585938"
$arr[7,7] = [double]"27922.25028991699"
$arr[7,8] = [double]"-0.02105680758348599"
$arr[7,9] = [double]"228.7657536548577"

$arr[8,0] = "2024-09-12"
$arr[8,2] = [double]"986.1500244140625"
$arr[8,3] = [double]"726.0499877929688"
$arr[8,4] = [double]"246.1499938964844"
$arr[8,5] = [double]"479.8500061035156"
$arr[8,6] = [double]"1403.150024414062"
$arr[8,7] = [double]"28182.30001831055"
$arr[8,8] = [double]"0.009313351384414074"
$arr[8,9] = [double]"230.8963295033658"

$arr[9,0] = "2024-09-13"
$arr[9,2] = [double]"992.0999755859375"
$arr[9,3] = [double]"724.25"
$arr[9,4] = [double]"245.6499938964844"
$arr[9,5] = [double]"485.3999938964844"
$arr[9,6] = [double]"1410.949951171875"
$arr[9,7] = [double]"28258.59951782227"
$arr[9,8] = [double]"0.002707355306775728"
$arr[9,9] = [double]"231.5214479063617"

$arr[10,0] = "2024-09-16"
$arr[10,2] = [double]"988.4000244140625"
$arr[10,3] = [double]"733.6500244140625"
$arr[10,4] = [double]"243.8000030517578"
$arr[10,5] = [double]"489.9500122070312"
$arr[10,6] = [double]"1404.550048828125"
$arr[10,7] = [double]"28313.45072937012"
$arr[10,8] = [double]"0.00194104493795801"
$arr[10,9] = [double]"231.9708414408491"

$arr[11,0] = "2024-09-17"
$arr[11,2] = [double]"974.9500122070312"
$arr[11,3] = [double]"745.4000244140625"
$arr[11,4] = [double]"240.8000030517578"
$arr[11,5] = [double]"482.2999877929688"
$arr[11,6] = [double]"1400.25"
$arr[11,7] = [double]"28196.30033874512"
$arr[11,8] = [double]"-0.004137623200533361"
$arr[11,9] = [double]"231.0110335054562"

$arr[12,0] = "2024-09-18"
$arr[12,2] = [double]"962.0499877929688"
$arr[12,3] = [double]"717.5499877929688"
$arr[12,4] = [double]"235.9499969482422"
$arr[12,5] = [double]"471.75"
$arr[12,6] = [double]"1391.300048828125"
$arr[12,7] = [double]"27572.89979553223"
$arr[12,8] = [double]"-0.02210930284198537"
$arr[12,9] = [double]"225.903540605844"

$arr[13,0] = "2024-09-19"
$arr[13,2] = [double]"967"
$arr[13,3] = [double]"728.5"
$arr[13,4] = [double]"237.5500030517578"
$arr[13,5] = [double]"459.9500122070312"
$arr[13,6] = [double]"1374.150024414062"
$arr[13,7] = [double]"27641.0502166748"
$arr[13,8] = [double]"0.002471645044516532"
$arr[13,9] = [double]"226.4618939725212"

$arr[14,0] = "2024-09-20"
$arr[14,2] = [double]"970.8499755859375"
$arr[14,3] = [double]"748.3499755859375"
$arr[14,4] = [double]"237.8500061035156"
$arr[14,5] = [double]"466.2999877929688"
$arr[14,6] = [double]"1380.550048828125"
$arr[14,7] = [double]"27960.69967651367"
$arr[14,8] = [double]"0.01156430227264067"
$arr[14,9] = [double]"229.0807677676541"

$arr[15,0] = "2024-09-23"
$arr[15,2] = [double]"971.7999877929688"
$arr[15,3] = [double]"750.2000122070312"
$arr[15,4] = [double]"236.4499969482422"
$arr[15,5] = [double]"471.1499938964844"
$arr[15,6] = [double]"1375.400024414062"
$arr[15,7] = [double]"27984.94996643066"
$arr[15,8] = [double]"0.0008672991090191444"
$arr[15,9] = [double]"229.2794493134324"

$arr[16,0] = "2024-09-24"
$arr[16,2] = [double]"977.2999877929688"
$arr[16,3] = [double]"735.9000244140625"
$arr[16,4] = [double]"237.3000030517578"
$arr[16,5] = [double]"476.7000122070312"
$arr[16,6] = [double]"1363.699951171875"
$arr[16,7] = [double]"27912.20024108887"
$arr[16,8] = [double]"-0.002599601765558408"
$arr[16,9] = [double]"228.683414052191"

$arr[17,0] = "2024-09-25"
$arr[17,2] = [double]"963.5999755859375"
$arr[17,3] = [double]"730.0499877929688"
$arr[17,4] = [double]"238.3500061035156"
$arr[17,5] = [double]"473.7000122070312"
$arr[17,6] = [double]"1365.400024414062"
$arr[17,7] = [double]"27741.79995727539"
$arr[17,8] = [double]"-0.006104867489544392"
$arr[17,9] = [double]"227.2873321123457"

$arr[18,0] = "2024-09-26"
$arr[18,2] = [double]"993.1500244140625"
$arr[18,3] = [double]"744.0999755859375"
$arr[18,4] = [double]"241.1999969482422"
$arr[18,5] = [double]"471.75"
$arr[18,6] = [double]"1329.949951171875"
$arr[18,7] = [double]"28111.79975891113"
$arr[18,8] = [double]"0.01333726730801793"
$arr[18,9] = [double]"230.3187240164543"

$arr[19,0] = "2024-09-27"
$arr[19,2] = [double]"993"
$arr[19,3] = [double]"735.4500122070312"
$arr[19,4] = [double]"239.5500030517578"
$arr[19,5] = [double]"497.2999877929688"
$arr[19,6] = [double]"1392.199951171875"
$arr[19,7] = [double]"28306.40000915527"
$arr[19,8] = [double]"0.006922368966521059"
$arr[19,9] = [double]"231.9130752039945"

# Ensure column A is treated as text so dates are not auto-converted to date serials
$ws.Range("A674:A693").NumberFormat = "@"
$ws.Range("A674:J693").Value() = $arr

# Re-apply the (unstyled) formatting used by the preceding data rows, then remove the
# now-empty column B cells that PasteSpecial would otherwise have created
$ws.Range("A673:J673").Copy()
$ws.Range("A674:J693").PasteSpecial(-4122)
$ws.Range("B674:B693").ClearContents()
$excel.CutCopyMode = 0